# error solve ifrs list
# Corrects the 인스코비 IFRS financial data rows (rows 2-9) in the active sheet.
# Rows 2-6: values are replaced with corrected figures (some cells cleared).
# Rows 7-9: all financial data cells (D:AJ) are cleared, keeping only A/B/C labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 258
$ws.Range("E2").Value = -58
$ws.Range("F2").Value = -58
$ws.Range("G2").Value = -79
$ws.Range("H2").Value = -80
$ws.Range("I2").Value = -80
$ws.Range("K2").Value = 239
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 209
$ws.Range("N2").Value = 209
$ws.Range("P2").Value = 298
$ws.Range("Q2").Value = -83
$ws.Range("R2").Value = -20
$ws.Range("S2").Value = 53
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = -83
$ws.Range("V2").Value = 7
$ws.Range("W2").Value = -22.6
$ws.Range("X2").Value = -30.88
$ws.Range("Y2").Value = -36.35
$ws.Range("Z2").Value = -31.8
$ws.Range("AA2").Value = 14.26
$ws.Range("AB2").Value = -32.36
$ws.Range("AC2").Value = -136
$ws.Range("AD2").Value = -14.54
$ws.Range("AE2").Value = 351
$ws.Range("AF2").Value = 5.65
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 59669884
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

$ws.Range("D3").Value = 163
$ws.Range("E3").Value = -48
$ws.Range("F3").Value = -48
$ws.Range("G3").Value = -245
$ws.Range("H3").Value = -245
$ws.Range("I3").Value = -235
$ws.Range("J3").Value = -10
$ws.Range("K3").Value = 486
$ws.Range("L3").Value = 233
$ws.Range("M3").Value = 254
$ws.Range("N3").Value = 259
$ws.Range("O3").Value = -5
$ws.Range("P3").Value = 377
$ws.Range("Q3").Value = -53
$ws.Range("R3").Value = -84
$ws.Range("S3").Value = 155
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = -54
$ws.Range("V3").Value = 146
$ws.Range("W3").Value = -29.12
$ws.Range("X3").Value = -150.45
$ws.Range("Y3").Value = -100.36
$ws.Range("Z3").Value = -67.65000000000001
$ws.Range("AA3").Value = 91.69
$ws.Range("AB3").Value = -31.99
$ws.Range("AC3").Value = -353
$ws.Range("AD3").Value = -4.46
$ws.Range("AE3").Value = 344
$ws.Range("AF3").Value = 4.58
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 75367969

$ws.Range("D4").Value = 371
$ws.Range("E4").Value = -25
$ws.Range("F4").Value = -25
$ws.Range("G4").Value = -63
$ws.Range("H4").Value = -69
$ws.Range("I4").Value = -61
$ws.Range("J4").Value = -7
$ws.Range("K4").Value = 452
$ws.Range("L4").Value = 133
$ws.Range("M4").Value = 319
$ws.Range("N4").Value = 332
$ws.Range("O4").Value = -13
$ws.Range("P4").Value = 424
$ws.Range("Q4").Value = -3
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = 17
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = -6
$ws.Range("V4").Value = 40
$ws.Range("W4").Value = -6.83
$ws.Range("X4").Value = -18.49
$ws.Range("Y4").Value = -20.69
$ws.Range("Z4").Value = -14.61
$ws.Range("AA4").Value = 41.55
$ws.Range("AB4").Value = -22.45
$ws.Range("AC4").Value = -77
$ws.Range("AD4").Value = -22.05
$ws.Range("AE4").Value = 391
$ws.Range("AF4").Value = 4.36
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 84826521

$ws.Range("D5").Value = 496
$ws.Range("E5").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = -5
$ws.Range("I5").Value = -6
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 459
$ws.Range("L5").Value = 121
$ws.Range("M5").Value = 337
$ws.Range("N5").Value = 346
$ws.Range("O5").Value = -8
$ws.Range("P5").Value = 436
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = -4
$ws.Range("S5").Value = -2
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 13
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = -0.4
$ws.Range("X5").Value = -0.92
$ws.Range("Y5").Value = -1.76
$ws.Range("Z5").Value = -1
$ws.Range("AA5").Value = 35.88
$ws.Range("AB5").Value = -21.69
$ws.Range("AC5").Value = -7
$ws.Range("AD5").Value = -327.67
$ws.Range("AE5").Value = 396
$ws.Range("AF5").Value = 5.68
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 87270467

$ws.Range("D6").Value = 586
$ws.Range("E6").Value = 53
$ws.Range("F6").Value = 53
$ws.Range("G6").Value = -76
$ws.Range("H6").Value = -81
$ws.Range("I6").Value = -83
$ws.Range("K6").Value = 574
$ws.Range("L6").Value = 266
$ws.Range("M6").Value = 309
$ws.Range("N6").Value = 309
$ws.Range("P6").Value = 449
$ws.Range("Q6").Value = 27
$ws.Range("R6").Value = -164
$ws.Range("S6").Value = 167
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 25
$ws.Range("V6").Value = 121
$ws.Range("W6").Value = 9.09
$ws.Range("X6").Value = -13.83
$ws.Range("Y6").Value = -25.43
$ws.Range("Z6").Value = -15.69
$ws.Range("AA6").Value = 85.94
$ws.Range("AB6").Value = -26.81
$ws.Range("AC6").Value = -94
$ws.Range("AD6").Value = -53.08
$ws.Range("AE6").Value = 348
$ws.Range("AF6").Value = 14.29
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 89723464

$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
